$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 8-10 (previously "MuSCs" sending-cluster block that is no longer present)
$ws.Range("A8:T10").Delete()

# Refresh rows 2-7 with updated TPM-derived values; row order now: FAPs block (rows 2-4), MuSCs block (rows 5-7)
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2660856666666667
$ws.Range("H2").Value = 0.798257
$ws.Range("I2").Value = 0.0824022153898117
$ws.Range("J2").Value = 0.1187122470109165
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1011536666666667
$ws.Range("N2").Value = 0.303461
$ws.Range("O2").Value = 0.007629860605400263
$ws.Range("P2").Value = 0.008254451482408482
$ws.Range("Q2").Value = 0.02691554083077777
$ws.Range("R2").Value = 0.242239867477
$ws.Range("S2").Value = 0.0006287174170004316
$ws.Range("T2").Value = 0.0009799044833193019

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2660856666666667
$ws.Range("H3").Value = 0.798257
$ws.Range("I3").Value = 0.0824022153898117
$ws.Range("J3").Value = 0.1187122470109165
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.14695433333333
$ws.Range("N3").Value = 30.440863
$ws.Range("O3").Value = 0.7653686681256785
$ws.Range("P3").Value = 0.8280227993585454
$ws.Range("Q3").Value = 2.699959108421222
$ws.Range("R3").Value = 24.299631975791
$ws.Range("S3").Value = 0.06306807384350546
$ws.Range("T3").Value = 0.09829644708812221

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2660856666666667
$ws.Range("H4").Value = 0.798257
$ws.Range("I4").Value = 0.0824022153898117
$ws.Range("J4").Value = 0.1187122470109165
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.0094955
$ws.Range("N4").Value = 6.018991
$ws.Range("O4").Value = 0.2270014712689213
$ws.Range("P4").Value = 0.1637227491590462
$ws.Range("Q4").Value = 0.8007836164478332
$ws.Range("R4").Value = 4.804701698686999
$ws.Range("S4").Value = 0.0187054241293058
$ws.Range("T4").Value = 0.01943589543947502

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.9630225
$ws.Range("H5").Value = 5.926045
$ws.Range("I5").Value = 0.9175977846101883
$ws.Range("J5").Value = 0.8812877529890836
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1011536666666667
$ws.Range("N5").Value = 0.303461
$ws.Range("O5").Value = 0.007629860605400263
$ws.Range("P5").Value = 0.008254451482408482
$ws.Range("Q5").Value = 0.2997205902908333
$ws.Range("R5").Value = 1.798323541745
$ws.Range("S5").Value = 0.007001143188399832
$ws.Range("T5").Value = 0.007274546999089181

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gdnf"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.9630225
$ws.Range("H6").Value = 5.926045
$ws.Range("I6").Value = 0.9175977846101883
$ws.Range("J6").Value = 0.8812877529890836
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 10.14695433333333
$ws.Range("N6").Value = 30.440863
$ws.Range("O6").Value = 0.7653686681256785
$ws.Range("P6").Value = 0.8280227993585454
$ws.Range("Q6").Value = 30.06565399613917
$ws.Range("R6").Value = 180.393923976835
$ws.Range("S6").Value = 0.702300594282173
$ws.Range("T6").Value = 0.7297263522704233

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gdnf"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.9630225
$ws.Range("H7").Value = 5.926045
$ws.Range("I7").Value = 0.9175977846101883
$ws.Range("J7").Value = 0.8812877529890836
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.0094955
$ws.Range("N7").Value = 6.018991
$ws.Range("O7").Value = 0.2270014712689213
$ws.Range("P7").Value = 0.1637227491590462
$ws.Range("Q7").Value = 8.91720288014875
$ws.Range("R7").Value = 35.668811520595
$ws.Range("S7").Value = 0.2082960471396154
$ws.Range("T7").Value = 0.1442868537195712
